# Update cryptos list - applies the diff of price/volume changes and the
# dogwifhat/Monero & EthereumClassic/PolygonEcosystemToken row swaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to be stored as text (matches source data which uses
    # formatted numeric strings like "153.30" or "64.801.40"), then restore
    # the default "Normal" style so no stray number-format style is left
    # attached to the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "64.801.40"
$ws.Range("E2").Value = "  +3.09%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.531.06"
$ws.Range("E3").Value = "  +2.71%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.01%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "580.37"
$ws.Range("E5").Value = "  +0.94%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "153.30"
$ws.Range("E6").Value = "  +3.81%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.541"
$ws.Range("E8").Value = "  +1.17%  "

# Row 9 - LidoStakedEther
Set-TextValue $ws.Range("D9") "2.533.93"
$ws.Range("E9").Value = "  +2.88%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.85%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -1.78%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.15%  "

# Row 13 - Cardano
$ws.Range("E13").Value = "  -0.27%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "29.28"
$ws.Range("E14").Value = "  +0.78%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +2.02%  "

# Row 16 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D16") "2.990.83"
$ws.Range("E16").Value = "  +2.70%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "64.207.06"
$ws.Range("E17").Value = "  +2.32%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.534.72"
$ws.Range("E18").Value = "  +2.84%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "7.98"
$ws.Range("E19").Value = "  +0.47%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.34%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  +3.58%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "330.03"
$ws.Range("E22").Value = "  +1.11%  "

# Row 23 - SuiNetwork
$ws.Range("E23").Value = "  +0.60%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  +0.07%  "

# Row 25 - Aptos
Set-TextValue $ws.Range("D25") "10.06"
$ws.Range("E25").Value = "  -0.08%  "

# Row 26 - Litecoin
Set-TextValue $ws.Range("D26") "65.71"
$ws.Range("E26").Value = "  +0.20%  "

# Row 27 - Bittensor
Set-TextValue $ws.Range("D27") "648.09"
$ws.Range("E27").Value = "  +0.42%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  +6.82%  "

# Row 29 - WrappedeETH
Set-TextValue $ws.Range("D29") "2.646.29"
$ws.Range("E29").Value = "  +2.45%  "

# Row 30 - Fetch.AI
$ws.Range("E30").Value = "  +5.06%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  -0.01%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +1.15%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +1.76%  "

# Row 34 - Kaspa
$ws.Range("E34").Value = "  +2.08%  "

# Row 35 - FirstDigitalUSD
Set-TextValue $ws.Range("D35") "0.998"
$ws.Range("E35").Value = "  -0.01%  "

# Row 36 - ImmutableX
$ws.Range("E36").Value = "  +1.99%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +2.20%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +3.52%  "

# Row 39 - becomes Monero (was dogwifhat)
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D39") "155.16"
$ws.Range("E39").Value = "  +1.97%  "

# Row 40 - becomes dogwifhat (was Monero)
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D40") "2.85"
$ws.Range("E40").Value = "  +1.46%  "

# Row 41 - becomes PolygonEcosystemToken (was EthereumClassic)
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D41") "0.373"
$ws.Range("E41").Value = "  +1.05%  "

# Row 42 - becomes EthereumClassic (was PolygonEcosystemToken)
$ws.Range("B42").Value = "EthereumClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D42") "18.99"
$ws.Range("E42").Value = "  +1.49%  "

# Row 43 - Stacks
$ws.Range("E43").Value = "  +3.60%  "

# Row 44 - Aave
Set-TextValue $ws.Range("D44") "163.76"
$ws.Range("E44").Value = "  +7.44%  "

# Row 45 - USDe
$ws.Range("E45").Value = "  +0.01%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -2.15%  "

# Row 47 - WhiteBITCoin
Set-TextValue $ws.Range("D47") "15.63"
$ws.Range("E47").Value = "  +2.28%  "

# Row 48 - Filecoin
Set-TextValue $ws.Range("D48") "3.66"
$ws.Range("E48").Value = "  +1.95%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "21.53"
$ws.Range("E49").Value = "  +4.97%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  +2.48%  "
